# Swap the contents of columns C and D (header + all data rows).
# In the original workbook:
#   Column C = codeforiati:group-code (e.g. "BE")
#   Column D = codeforiati:group-name (e.g. "Bélgica")
# After the edit:
#   Column C = codeforiati:group-name (e.g. "Bélgica")
#   Column D = codeforiati:group-code (e.g. "BE")
# i.e. the two columns simply exchange places, for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 1) { $lastRow = 1 }

$colC = $ws.Range("C1:C$lastRow")
$colD = $ws.Range("D1:D$lastRow")

$valuesC = $colC.Value2
$valuesD = $colD.Value2

$colC.Value2 = $valuesD
$colD.Value2 = $valuesC
